$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.816.46"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").Value = "'3.428.27"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'574.25"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").Value = "'158.64"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").Value = "  +2.75%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "'3.428.41"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "'4.019.26"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "'0.0000188"
$ws.Range("E15").Value = "  -3.74%  "
$ws.Range("D16").Value = "'27.74"
$ws.Range("E16").Value = "  -3.98%  "
$ws.Range("D17").Value = "'64.809.53"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").Value = "'3.422.91"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "'13.85"
$ws.Range("E20").Value = "  -3.22%  "
$ws.Range("D21").Value = "'381.15"
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("D22").Value = "'8.01"
$ws.Range("E22").Value = "  -2.90%  "
$ws.Range("D23").Value = "'0.549"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "'72.32"
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("E26").Value = "  -5.27%  "
$ws.Range("D27").Value = "'10.06"
$ws.Range("E27").Value = "  +3.38%  "
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +2.81%  "
$ws.Range("D31").Value = "'6.18"
$ws.Range("E31").Value = "  -4.11%  "
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("D33").Value = "'23.27"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("D34").Value = "'7.09"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").Value = "'1.60"
$ws.Range("E35").Value = "  +2.92%  "
$ws.Range("D36").Value = "'160.58"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("D38").Value = "'2.927.31"
$ws.Range("E38").Value = "  -5.02%  "
$ws.Range("D39").Value = "'0.0757"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("D40").Value = "'6.77"
$ws.Range("E40").Value = "  +4.18%  "
$ws.Range("D41").Value = "'26.42"
$ws.Range("E41").Value = "  -3.29%  "
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("D43").Value = "'43.01"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("D45").Value = "'26.01"
$ws.Range("E45").Value = "  +1.29%  "
$ws.Range("D46").Value = "'0.772"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").Value = "'319.42"
$ws.Range("E47").Value = "  +2.55%  "
$ws.Range("D48").Value = "'2.24"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("D49").Value = "'1.07"
$ws.Range("E49").Value = "  -4.59%  "
$ws.Range("D50").Value = "'0.108"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("E51").Value = "  -2.19%  "
